$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new shared-string notes for the new rows
$ws.Range("D14").Value = "Working on logging function. Running 100 simulations with random behaviour to generate log of base case."
$ws.Range("D15").Value = "Created a spreadsheet to analyse baseline data, found a bug in simulator from the statistical data. Fixed simulator and generated another data set."

# Row 14: 2017-06-19 (serial 42905)
$ws.Range("A14").Value = 42905
$ws.Range("B14").Formula = "=A14"
$ws.Range("C14").Value = 3

# Row 15: 2017-06-20 (serial 42906)
$ws.Range("A15").Value = 42906
$ws.Range("B15").Formula = "=A15"
$ws.Range("C15").Value = 1

# Copy formatting/styles from row 13 (existing last row) down to the new rows
$ws.Range("A13:D13").Copy() | Out-Null
$ws.Range("A14:D14").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("A13:D13").Copy() | Out-Null
$ws.Range("A15:D15").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Row heights for the new wrap-text note rows (matches the 28.5pt used by
# other two-line note rows such as row 11)
$ws.Rows.Item(14).RowHeight = 28.5
$ws.Rows.Item(15).RowHeight = 28.5

# Update selection to F2
$ws.Range("F2").Select() | Out-Null

$wb.Save()
